$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.549.36'
$ws.Range("E2").Value = '  -0.11%  '

$ws.Range("D3").Value = '3.675.73'
$ws.Range("E3").Value = '  -0.84%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''620.69'
$ws.Range("E5").Value = '  -7.82%  '

$ws.Range("D6").Value = '''159.46'
$ws.Range("E6").Value = '  -1.58%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  -0.25%  '

$ws.Range("D9").Value = '''0.145'
$ws.Range("E9").Value = '  -1.35%  '

$ws.Range("D10").Value = '''7.21'
$ws.Range("E10").Value = '  +1.41%  '

$ws.Range("D11").Value = '''0.441'
$ws.Range("E11").Value = '  -0.69%  '

$ws.Range("E12").Value = '  -2.47%  '

$ws.Range("D13").Value = '4.294.56'
$ws.Range("E13").Value = '  -0.88%  '

$ws.Range("D14").Value = '''32.47'
$ws.Range("E14").Value = '  -1.21%  '

$ws.Range("D15").Value = '3.667.22'
$ws.Range("E15").Value = '  -0.96%  '

$ws.Range("D16").Value = '69.559.78'
$ws.Range("E16").Value = '  -0.14%  '

$ws.Range("E17").Value = '  +0.68%  '

$ws.Range("D18").Value = '''6.52'
$ws.Range("E18").Value = '  -0.23%  '

$ws.Range("D19").Value = '''15.92'
$ws.Range("E19").Value = '  -2.55%  '

$ws.Range("D20").Value = '''10.32'
$ws.Range("E20").Value = '  +5.11%  '

$ws.Range("D21").Value = '''469.31'
$ws.Range("E21").Value = '  -1.03%  '

$ws.Range("D22").Value = '''0.651'
$ws.Range("E22").Value = '  -0.60%  '

$ws.Range("D23").Value = '''79.59'
$ws.Range("E23").Value = '  -1.08%  '

$ws.Range("D24").Value = '3.821.16'
$ws.Range("E24").Value = '  -0.91%  '

$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").Value = '''11.14'
$ws.Range("E26").Value = '  +1.09%  '

$ws.Range("D27").Value = '''0.0000123'
$ws.Range("E27").Value = '  -4.04%  '

$ws.Range("D28").Value = '''8.68'
$ws.Range("E28").Value = '  -5.29%  '

$ws.Range("E29").Value = '  -3.06%  '

$ws.Range("E30").Value = '  -4.08%  '

$ws.Range("E31").Value = '  +0.18%  '

$ws.Range("E32").Value = '  -2.19%  '

$ws.Range("D33").Value = '''26.64'
$ws.Range("E33").Value = '  -1.19%  '

$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = '''6.41'
$ws.Range("E34").Value = '  -2.83%  '

$ws.Range("B35").Value = 'RenzoRestakedETH'
$ws.Range("C35").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D35").Value = '3.674.52'
$ws.Range("E35").Value = '  -0.62%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '''0.162'
$ws.Range("E36").Value = '  -3.20%  '

$ws.Range("D37").Value = '''8.30'
$ws.Range("E37").Value = '  -2.94%  '

$ws.Range("E38").Value = '  -0.02%  '

$ws.Range("D39").Value = '''178.75'
$ws.Range("E39").Value = '  +2.73%  '

$ws.Range("E41").Value = '  -1.46%  '

$ws.Range("D42").Value = '''5.81'
$ws.Range("E42").Value = '  -5.25%  '

$ws.Range("D43").Value = '''0.0896'
$ws.Range("E43").Value = '  -1.91%  '

$ws.Range("E44").Value = '  -1.61%  '

$ws.Range("D45").Value = '''29.40'
$ws.Range("E45").Value = '  +6.08%  '

$ws.Range("D46").Value = '''46.76'
$ws.Range("E46").Value = '  -0.66%  '

$ws.Range("D47").Value = '''2.72'
$ws.Range("E47").Value = '  -1.63%  '

$ws.Range("D48").Value = '''7.87'
$ws.Range("E48").Value = '  -0.16%  '

$ws.Range("D49").Value = '''0.000266'
$ws.Range("E49").Value = '  -5.30%  '

$ws.Range("E50").Value = '  -4.82%  '

$ws.Range("E51").Value = '  -6.72%  '
